$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.Value = "'" + $val
    $c.Style = "Normal"
}

Set-TextValue "D2" "57.957.93"
Set-TextValue "E2" "  -0.03%  "
Set-TextValue "D3" "2.350.17"
Set-TextValue "E3" "  +1.01%  "
Set-TextValue "E4" "  -0.27%  "
Set-TextValue "D5" "541.30"
Set-TextValue "E5" "  -0.14%  "
Set-TextValue "D6" "134.72"
Set-TextValue "E6" "  -0.11%  "
Set-TextValue "E7" "  +0.22%  "
Set-TextValue "D8" "0.568"
Set-TextValue "E9" "  +0.54%  "
Set-TextValue "D10" "5.53"
Set-TextValue "E10" "  +3.07%  "
Set-TextValue "E11" "  -1.71%  "
Set-TextValue "E12" "  +0.79%  "
Set-TextValue "D13" "23.79"
Set-TextValue "E13" "  +1.25%  "
Set-TextValue "D14" "2.770.10"
Set-TextValue "E14" "  +0.72%  "
Set-TextValue "D15" "57.895.73"
Set-TextValue "E15" "  +0.05%  "
Set-TextValue "E16" "  +0.93%  "
Set-TextValue "D17" "2.371.28"
Set-TextValue "E17" "  +0.68%  "
Set-TextValue "E18" "  +1.63%  "
Set-TextValue "E19" "  +2.32%  "
Set-TextValue "D20" "329.83"
Set-TextValue "E20" "  -2.52%  "
Set-TextValue "D21" "6.73"
Set-TextValue "E21" "  -1.23%  "
Set-TextValue "E22" "  +0.48%  "
Set-TextValue "D23" "62.78"
Set-TextValue "E23" "  +0.95%  "
Set-TextValue "E24" "  -2.60%  "
Set-TextValue "D25" "1.00"
Set-TextValue "E25" "  +0.15%  "
Set-TextValue "E26" "  -1.44%  "
Set-TextValue "D27" "1.35"
Set-TextValue "E27" "  -5.17%  "
Set-TextValue "E28" "  +0.38%  "
Set-TextValue "D29" "170.01"
Set-TextValue "E29" "  -0.57%  "
Set-TextValue "E30" "  +0.13%  "
Set-TextValue "E31" "  -0.53%  "
Set-TextValue "E32" "  +1.17%  "
Set-TextValue "E33" "  -0.98%  "
Set-TextValue "E34" "  +0.06%  "
Set-TextValue "E35" "  +0.33%  "
Set-TextValue "D36" "4.20"
Set-TextValue "E36" "  +1.41%  "
Set-TextValue "E37" "  -1.50%  "
Set-TextValue "E38" "  +0.20%  "
Set-TextValue "D39" "39.09"
Set-TextValue "D40" "142.64"
Set-TextValue "E40" "  -3.96%  "
Set-TextValue "B41" "PolygonEcosystemToken"
Set-TextValue "C41" "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
Set-TextValue "D41" "0.378"
Set-TextValue "E41" "  -0.18%  "
Set-TextValue "B42" "Filecoin"
Set-TextValue "C42" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D42" "3.65"
Set-TextValue "E42" "  +0.64%  "
Set-TextValue "D43" "288.88"
Set-TextValue "E43" "  +1.77%  "
Set-TextValue "D44" "0.0949"
Set-TextValue "E44" "  +1.73%  "
Set-TextValue "E45" "  +0.78%  "
Set-TextValue "D46" "19.12"
Set-TextValue "E46" "  -0.18%  "
Set-TextValue "E47" "  +1.55%  "
Set-TextValue "D48" "0.0222"
Set-TextValue "E48" "  +1.58%  "
Set-TextValue "D49" "0.381"
Set-TextValue "E49" "  -0.22%  "
Set-TextValue "E50" "  +0.53%  "
Set-TextValue "E51" "  -0.62%  "
